$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 917
$ws.Range("I21").Value = 917
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 917
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -449

$ws.Range("H23").Value = 917
$ws.Range("I23").Value = 917
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 917
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -683

$ws.Range("H34").Value = 1636982
$ws.Range("I34").Value = 1636982
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1636982
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1636779
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 1636982
$ws.Range("I36").Value = 1636982
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1636982
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1636267
$ws.Range("N36").ClearContents()

$ws.Range("H88").Value = 1725
$ws.Range("I88").Value = 1300
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 1300
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -894
$ws.Range("N88").Value = -3812

$ws.Range("H91").Value = 1725
$ws.Range("I91").Value = 1300
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 1300
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = 104
$ws.Range("N91").Value = -5808

$ws.Range("H111").Value = 3875
$ws.Range("I111").Value = 3833.3333
$ws.Range("J111").Value = 4000
$ws.Range("K111").Value = 11499.9999
$ws.Range("L111").Value = 12000
$ws.Range("M111").Value = -8432.999899999999
$ws.Range("N111").Value = -18134

$ws.Range("H113").Value = 2949.25
$ws.Range("I113").Value = 1970
$ws.Range("J113").Value = 3014.5334
$ws.Range("K113").Value = 1970
$ws.Range("L113").Value = 3014.5334
$ws.Range("M113").Value = 1284
$ws.Range("N113").Value = -9522.5334

$ws.Range("H116").Value = 2124.75
$ws.Range("I116").Value = 1908
$ws.Range("J116").Value = 2775
$ws.Range("K116").Value = 1908
$ws.Range("L116").Value = 2775
$ws.Range("M116").Value = 1534
$ws.Range("N116").Value = -9659

$ws.Range("H125").Value = 11597.667
$ws.Range("I125").Value = 478.6
$ws.Range("J125").Value = 25496.5
$ws.Range("K125").Value = 4307.400000000001
$ws.Range("L125").Value = 229468.5
$ws.Range("M125").Value = -1847.400000000001
$ws.Range("N125").Value = -234388.5

$ws.Range("H127").Value = 71429370
$ws.Range("I127").Value = 613
$ws.Range("J127").Value = 142858130
$ws.Range("K127").Value = 1839
$ws.Range("L127").Value = 428574390
$ws.Range("M127").Value = 3121
$ws.Range("N127").Value = -428584310

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1143.5
$ws.Range("I110").Value = 1104.1177
$ws.Range("J110").Value = 1366.6666
$ws.Range("K110").Value = 1104.1177
$ws.Range("L110").Value = 1366.6666
$ws.Range("M110").Value = 940.8823
$ws.Range("N110").Value = -5456.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2192.0527
$ws.Range("I16").Value = 1689.7
$ws.Range("J16").Value = 2750.2222
$ws.Range("K16").Value = 1689.7
$ws.Range("L16").Value = 2750.2222
$ws.Range("M16").Value = -1402.7
$ws.Range("N16").Value = -3324.2222

$ws.Range("H22").Value = 209.625
$ws.Range("I22").Value = 242.5
$ws.Range("J22").Value = 176.75
$ws.Range("K22").Value = 242.5
$ws.Range("L22").Value = 176.75
$ws.Range("M22").Value = 107.5
$ws.Range("N22").Value = -876.75

$ws.Range("H99").Value = 3690
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 3690
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 3690
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -6686

$ws.Range("H105").Value = 722.8570999999999
$ws.Range("I105").Value = 678.3333
$ws.Range("J105").Value = 990
$ws.Range("K105").Value = 678.3333
$ws.Range("L105").Value = 990
$ws.Range("M105").Value = 1068.6667
$ws.Range("N105").Value = -4484

$ws.Range("H113").Value = 2192.0527
$ws.Range("I113").Value = 1689.7
$ws.Range("J113").Value = 2750.2222
$ws.Range("K113").Value = 1689.7
$ws.Range("L113").Value = 2750.2222
$ws.Range("M113").Value = 480.3
$ws.Range("N113").Value = -7090.2222

$ws.Range("H126").Value = 3690
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3690
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 11070
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -16010

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4009.6
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4009.6
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 12028.8
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -14024.8

$ws.Range("H78").Value = 4009.6
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4009.6
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 36086.4
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -46070.4

$ws.Range("H112").Value = 2392.7778
$ws.Range("I112").Value = 1218.6666
$ws.Range("J112").Value = 2979.8333
$ws.Range("K112").Value = 3655.9998
$ws.Range("L112").Value = 8939.499899999999
$ws.Range("M112").Value = -2547.9998
$ws.Range("N112").Value = -11155.4999

$ws.Range("H115").Value = 2833.3333
$ws.Range("I115").Value = 500
$ws.Range("J115").Value = 4000
$ws.Range("K115").Value = 1500
$ws.Range("L115").Value = 12000
$ws.Range("M115").Value = -325
$ws.Range("N115").Value = -14350

$ws.Range("H121").Value = 889.875
$ws.Range("I121").Value = 682.5
$ws.Range("J121").Value = 959
$ws.Range("K121").Value = 2047.5
$ws.Range("L121").Value = 2877
$ws.Range("M121").Value = -737.5
$ws.Range("N121").Value = -5497

$ws.Range("H123").Value = 2740.6667
$ws.Range("I123").Value = 1000
$ws.Range("J123").Value = 2958.25
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 8874.75
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -13774.75

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("M128").ClearContents()

$ws.Range("H139").Value = 2382.087
$ws.Range("I139").Value = 1519.3334
$ws.Range("J139").Value = 3999.75
$ws.Range("K139").Value = 4558.0002
$ws.Range("L139").Value = 11999.25
$ws.Range("M139").Value = 581.9997999999996
$ws.Range("N139").Value = -22279.25

$ws.Range("H141").Value = 3490
$ws.Range("I141").Value = 3490
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10470
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5290
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2641.8235
$ws.Range("I113").Value = 2739.3076
$ws.Range("J113").Value = 2325
$ws.Range("K113").Value = 2739.3076
$ws.Range("L113").Value = 2325
$ws.Range("M113").Value = -569.3076000000001
$ws.Range("N113").Value = -6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 609.5
$ws.Range("I22").Value = 278.75
$ws.Range("J22").Value = 774.875
$ws.Range("K22").Value = 278.75
$ws.Range("L22").Value = 774.875
$ws.Range("M22").Value = 16.25
$ws.Range("N22").Value = -1364.875

$ws.Range("H27").Value = 609.5
$ws.Range("I27").Value = 278.75
$ws.Range("J27").Value = 774.875
$ws.Range("K27").Value = 278.75
$ws.Range("L27").Value = 774.875
$ws.Range("M27").Value = -171.75
$ws.Range("N27").Value = -988.875

$ws.Range("H61").Value = 1755523.4
$ws.Range("I61").Value = 1652491.9
$ws.Range("J61").Value = 2002799
$ws.Range("K61").Value = 1652491.9
$ws.Range("L61").Value = 2002799
$ws.Range("M61").Value = -1652289.9
$ws.Range("N61").Value = -2003203

$ws.Range("H113").Value = 1755523.4
$ws.Range("I113").Value = 1652491.9
$ws.Range("J113").Value = 2002799
$ws.Range("K113").Value = 1652491.9
$ws.Range("L113").Value = 2002799
$ws.Range("M113").Value = -1650321.9
$ws.Range("N113").Value = -2007139

$ws.Range("H133").Value = 25868
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 25868
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 25868
$ws.Range("N133").Value = -30928

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 639.70966
$ws.Range("I113").Value = 438.27274
$ws.Range("J113").Value = 750.5
$ws.Range("K113").Value = 1314.81822
$ws.Range("L113").Value = 2251.5
$ws.Range("M113").Value = 855.1817799999999
$ws.Range("N113").Value = -6591.5
